# Atualização automática de pedidos - 30/05/2025 09:06
# Adds REQ-013 to both the "Pedidos" and "Itens" sheets, and normalizes the
# previous REQ-012 RACK/seccao value (row 13) from text to a real number on
# both sheets.

$wb = $excel.ActiveWorkbook
$wsPedidos = $wb.Worksheets.Item("Pedidos")
$wsItens   = $wb.Worksheets.Item("Itens")

# --- Pedidos: fix D13 (RACK) so it is stored as a number, not text -------
$wsPedidos.Cells.Item(13, 4).Value = 12

# --- Pedidos: append REQ-013 as row 14 ------------------------------------
$wsPedidos.Cells.Item(14, 1).Value = "REQ-013"
$wsPedidos.Cells.Item(14, 2).Value = "30/05/2025 09:06"
$wsPedidos.Cells.Item(14, 3).Value = "Renault"
# RACK stays textual ("12"), same quirk as the previous row before the fix.
# The leading apostrophe forces text storage; resetting the style afterwards
# keeps the cell on the sheet's default (un-styled) look, matching the rest
# of the column.
$wsPedidos.Cells.Item(14, 4).Value = "'12"
$wsPedidos.Cells.Item(14, 4).Style = "Normal"
$wsPedidos.Cells.Item(14, 5).Value = "R12-LA-A1"
$wsPedidos.Cells.Item(14, 6).Value = "washington vieira"
$wsPedidos.Cells.Item(14, 7).Value = ""
$wsPedidos.Cells.Item(14, 8).Value = "Pendente"
$wsPedidos.Cells.Item(14, 9).Value = ""
$wsPedidos.Cells.Item(14, 10).Value = ""

# --- Itens: fix D13 (seccao) so it is stored as a number, not text -------
$wsItens.Cells.Item(13, 4).Value = 0.5

# --- Itens: append REQ-013 as row 14 --------------------------------------
$wsItens.Cells.Item(14, 1).Value = "REQ-013"
$wsItens.Cells.Item(14, 2).Value = "ACOR2Z-0.35-GY"
$wsItens.Cells.Item(14, 3).Value = "180EX606941"
$wsItens.Cells.Item(14, 4).Value = "'0.35"
$wsItens.Cells.Item(14, 4).Style = "Normal"
$wsItens.Cells.Item(14, 5).Value = "GY"
$wsItens.Cells.Item(14, 6).Value = 1
